$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue 'D2' '61.385.71'
Set-TextValue 'D3' '2.391.35'
$ws.Range('E3').Value = '  -3.89%  '
$ws.Range('E4').Value = '  -0.07%  '
Set-TextValue 'D5' '549.99'
$ws.Range('E5').Value = '  -1.20%  '
Set-TextValue 'D6' '142.02'
$ws.Range('E6').Value = '  -4.55%  '
$ws.Range('E7').Value = '  -0.07%  '
Set-TextValue 'D8' '0.533'
$ws.Range('E8').Value = '  -11.37%  '
Set-TextValue 'D9' '2.388.05'
$ws.Range('E9').Value = '  -3.98%  '
$ws.Range('E10').Value = '  -3.04%  '
$ws.Range('E11').Value = '  +0.21%  '
$ws.Range('E12').Value = '  -3.50%  '
$ws.Range('E13').Value = '  -3.28%  '
Set-TextValue 'D14' '25.54'
$ws.Range('E14').Value = '  -3.73%  '
Set-TextValue 'D15' '2.820.45'
$ws.Range('E15').Value = '  -3.95%  '
$ws.Range('E16').Value = '  -2.63%  '
Set-TextValue 'D17' '60.829.78'
$ws.Range('E17').Value = '  -1.31%  '
Set-TextValue 'D18' '2.390.65'
$ws.Range('E18').Value = '  -3.77%  '
Set-TextValue 'D19' '10.76'
$ws.Range('E19').Value = '  -4.73%  '
$ws.Range('E20').Value = '  -2.61%  '
Set-TextValue 'D21' '319.12'
$ws.Range('E21').Value = '  -1.37%  '
$ws.Range('E22').Value = '  -6.16%  '
$ws.Range('E23').Value = '  +0.01%  '
Set-TextValue 'D24' '1.91'
$ws.Range('E24').Value = '  -0.54%  '
Set-TextValue 'D25' '63.86'
$ws.Range('E25').Value = '  -0.85%  '
Set-TextValue 'D26' '8.24'
$ws.Range('E26').Value = '  +3.93%  '
Set-TextValue 'D27' '0.999'
$ws.Range('E27').Value = '  +0.13%  '
Set-TextValue 'D28' '2.508.41'
$ws.Range('E28').Value = '  -3.59%  '
$ws.Range('B29').Value = 'PEPE'
$ws.Range('C29').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue 'D29' '0.0₃0932'
$ws.Range('E29').Value = '  -8.84%  '
$ws.Range('B30').Value = 'Bittensor'
$ws.Range('C30').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue 'D30' '530.68'
$ws.Range('E30').Value = '  -4.62%  '
$ws.Range('B31').Value = 'Fetch.AI'
$ws.Range('C31').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue 'D31' '1.44'
$ws.Range('E31').Value = '  -6.11%  '
$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue 'D32' '8.12'
$ws.Range('E32').Value = '  -3.30%  '
$ws.Range('E33').Value = '  -4.36%  '
Set-TextValue 'D34' '1.85'
$ws.Range('E34').Value = '  -4.11%  '
$ws.Range('E35').Value = '  -1.22%  '
$ws.Range('E36').Value = '  -0.12%  '
$ws.Range('E37').Value = '  -7.93%  '
$ws.Range('E38').Value = '  -5.51%  '
Set-TextValue 'D39' '0.376'
$ws.Range('E39').Value = '  -2.78%  '
Set-TextValue 'D40' '1.84'
$ws.Range('E40').Value = '  +3.44%  '
Set-TextValue 'D41' '18.14'
$ws.Range('E41').Value = '  -2.91%  '
Set-TextValue 'D42' '140.14'
$ws.Range('E42').Value = '  -4.65%  '
$ws.Range('E43').Value = '  +0.01%  '
Set-TextValue 'D44' '40.38'
$ws.Range('E44').Value = '  -0.66%  '
$ws.Range('B45').Value = 'dogwifhat'
$ws.Range('C45').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue 'D45' '2.15'
$ws.Range('E45').Value = '  -11.87%  '
$ws.Range('B46').Value = 'Filecoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 'D46' '3.63'
$ws.Range('E46').Value = '  -1.26%  '
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue 'D47' '141.00'
$ws.Range('E47').Value = '  -5.23%  '
Set-TextValue 'D48' '20.18'
$ws.Range('E48').Value = '  -9.32%  '
Set-TextValue 'D49' '0.0519'
$ws.Range('E49').Value = '  -4.92%  '
$ws.Range('E50').Value = '  -3.77%  '
$ws.Range('E51').Value = '  -4.36%  '
